# Edit: Exclude international aviation from CB constraint and net-zero constraint
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "config" sheet: insert two new rows (new row 8, new row 13) and populate
#    them. Writing "UC_CUMFLO" before "Pset_PN"/"T-A*INT*" fixes shared
#    string indices at 82/83/84 to match the target workbook.
# ---------------------------------------------------------------------------
$cfg = $wb.Worksheets.Item("config")
$cfg.Rows("8:8").Insert(-4121)
$cfg.Rows("13:13").Insert(-4121)

$cfg.Range("B13").Value2 = "UC_CUMFLO"
$cfg.Range("C13").Value2 = -1

$cfg.Range("B8").Value2 = "Pset_PN"
$cfg.Range("C8").Value2 = "T-A*INT*"

# ---------------------------------------------------------------------------
# helper to rebuild a scenario sheet ("single" / "multi")
# ---------------------------------------------------------------------------
function Update-ScenarioSheet {
    param($ws, $suffix, $i7UsesAbsRow5)

    # Insert the new Pset_PN column (E) and the new UC_CUMFLO column (J)
    $ws.Columns("E:E").Insert(-4161)
    $ws.Columns("J:J").Insert(-4161)

    # Row 5 headers
    $ws.Range("E5").Value2 = "Pset_PN"
    $ws.Range("J5").Value2 = "UC_CUMFLO"

    # Row 7 formulas (existing "Single value" row), now pointing at $C$14
    $ws.Range("B7").Formula = "=VLOOKUP(B5, config!`$B`$4:`$C`$14,2,FALSE) & ""_$suffix"""
    $ws.Range("C7").Formula = "=VLOOKUP(C`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("D7").Formula = "=VLOOKUP(D`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("F7").Formula = "=VLOOKUP(F`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("G7").Formula = "=VLOOKUP(G5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("H7").Formula = "=VLOOKUP(H5, config!`$B`$4:`$C`$14,2,FALSE)"
    if ($i7UsesAbsRow5) {
        $ws.Range("I7").Formula = "=VLOOKUP(I`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    } else {
        $ws.Range("I7").Formula = "=VLOOKUP(I5, config!`$B`$4:`$C`$14,2,FALSE)"
    }
    $ws.Range("E7").ClearFormats()
    $ws.Range("K7").Formula = "=VLOOKUP(""Value"", config!`$B`$4:`$C`$14,2,FALSE)*1000"
    $ws.Range("J7").ClearFormats()
    $ws.Range("L7").Formula = "=VLOOKUP(L5, config!`$B`$4:`$C`$14,2,FALSE) & "" - $suffix"""

    # Row 8 (new "T-A*INT*" / UC_CUMFLO = -1 row)
    $ws.Range("C8").Formula = "=VLOOKUP(C`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("D8").Formula = "=VLOOKUP(D`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("E8").Formula = "=VLOOKUP(E`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("F8").Formula = "=VLOOKUP(F`$5, config!`$B`$4:`$C`$14,2,FALSE)"
    $ws.Range("J8").Formula = "=VLOOKUP(J`$5, config!`$B`$4:`$C`$14,2,FALSE)"
}

# ---------------------------------------------------------------------------
# 2. "single" sheet
# ---------------------------------------------------------------------------
$single = $wb.Worksheets.Item("single")
Update-ScenarioSheet -ws $single -suffix "Single" -i7UsesAbsRow5 $true
$single.Range("G26").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. "multi" sheet
# ---------------------------------------------------------------------------
$multi = $wb.Worksheets.Item("multi")
Update-ScenarioSheet -ws $multi -suffix "Multi" -i7UsesAbsRow5 $false
$multi.Range("M27").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Make "multi" the active / selected sheet (moves tabSelected + activeTab)
# ---------------------------------------------------------------------------
$multi.Activate()
